$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (old rows 13 "Motor function" and 14 "CD39-Diplotypes")
# since the data set shrinks from 14 rows to 12 rows and values shift up.
$ws.Rows("13:14").Delete()

# Update labels and VIF values for rows 2-12 to the new dataset.
$ws.Range("A2").Value = "Sex"
$ws.Range("B2").Value = 3.359163680525868

$ws.Range("A3").Value = "Age"
$ws.Range("B3").Value = 5.996217456570353

$ws.Range("A4").Value = "Motor function"
$ws.Range("B4").Value = 2.204784773044462

$ws.Range("A5").Value = "Dermatological symptoms"
$ws.Range("B5").Value = 1.286065027082814

$ws.Range("A6").Value = "Arthralgia"
$ws.Range("B6").Value = 1.359287313735909

$ws.Range("A7").Value = "Urinary function impairment"
$ws.Range("B7").Value = 3.28598490754178

$ws.Range("A8").Value = "Lower limb pain"
$ws.Range("B8").Value = 1.540918493680095

$ws.Range("A9").Value = "Paresthesia"
$ws.Range("B9").Value = 1.321305717202734

$ws.Range("A10").Value = "Lower back spine pain"
$ws.Range("B10").Value = 1.818504763135346

$ws.Range("A11").Value = "Limbs paresis or weakness"
$ws.Range("B11").Value = 2.216266773874763

$ws.Range("A12").Value = "CD39-Diplotypes"
$ws.Range("B12").Value = 6.38178132333618
